# Auto-generated script: apply scheduled-runner market data updates
# to the Typhon_Profits sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 314.0909
$ws.Range("I5").Value = 38
$ws.Range("K5").Value = 38
$ws.Range("M5").Value = 77
$ws.Range("H11").Value = 95.14286
$ws.Range("I11").Value = 95.14286
$ws.Range("K11").Value = 95.14286
$ws.Range("M11").Value = 44.85714
$ws.Range("H17").Value = 475.63635
$ws.Range("J17").Value = 472.16666
$ws.Range("L17").Value = 1416.49998
$ws.Range("N17").Value = -1752.49998
$ws.Range("H98").Value = 778.8889
$ws.Range("I98").Value = 778.8889
$ws.Range("K98").Value = 778.8889
$ws.Range("M98").Value = 719.1111
$ws.Range("H116").Value = 4159.385
$ws.Range("I116").Value = 1652
$ws.Range("J116").Value = 5726.5
$ws.Range("K116").Value = 1652
$ws.Range("L116").Value = 5726.5
$ws.Range("M116").Value = 1790
$ws.Range("N116").Value = -12610.5
$ws.Range("H122").Value = 778.8889
$ws.Range("I122").Value = 778.8889
$ws.Range("K122").Value = 2336.6667
$ws.Range("M122").Value = 113.3332999999998
$ws.Range("H129").Value = 833.6923
$ws.Range("H138").Value = 3493.4546
$ws.Range("I138").Value = 1546.75
$ws.Range("J138").Value = 3926.0557
$ws.Range("K138").Value = 4640.25
$ws.Range("L138").Value = 11778.1671
$ws.Range("M138").Value = 499.75
$ws.Range("N138").Value = -22058.1671

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1154.1765
$ws.Range("I2").Value = 909.12823
$ws.Range("J2").Value = 1950.5834
$ws.Range("K2").Value = 909.12823
$ws.Range("L2").Value = 1950.5834
$ws.Range("M2").Value = -796.12823
$ws.Range("N2").Value = -2176.5834
$ws.Range("H32").Value = 17268.295
$ws.Range("I32").Value = 20764.64
$ws.Range("K32").Value = 20764.64
$ws.Range("M32").Value = -20477.64
$ws.Range("H63").Value = 3474411
$ws.Range("I63").Value = 2385.4285
$ws.Range("J63").Value = 15626500
$ws.Range("K63").Value = 2385.4285
$ws.Range("L63").Value = 15626500
$ws.Range("M63").Value = -1699.4285
$ws.Range("N63").Value = -15627872
$ws.Range("H66").Value = 3474411
$ws.Range("I66").Value = 2385.4285
$ws.Range("J66").Value = 15626500
$ws.Range("K66").Value = 11927.1425
$ws.Range("L66").Value = 78132500
$ws.Range("M66").Value = -8495.1425
$ws.Range("N66").Value = -78139364
$ws.Range("H116").Value = 1154.1765
$ws.Range("I116").Value = 909.12823
$ws.Range("J116").Value = 1950.5834
$ws.Range("K116").Value = 909.12823
$ws.Range("L116").Value = 1950.5834
$ws.Range("M116").Value = 1384.87177
$ws.Range("N116").Value = -6538.5834

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1154.1765
$ws.Range("I3").Value = 909.12823
$ws.Range("J3").Value = 1950.5834
$ws.Range("K3").Value = 909.12823
$ws.Range("L3").Value = 1950.5834
$ws.Range("M3").Value = -795.12823
$ws.Range("N3").Value = -2178.5834
$ws.Range("H99").Value = 2340
$ws.Range("I99").Value = 2250
$ws.Range("K99").Value = 2250
$ws.Range("M99").Value = -752
$ws.Range("H107").Value = 750
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 1000
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 920
$ws.Range("N107").Value = -4340

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11931.978
$ws.Range("I31").Value = 22270.65
$ws.Range("K31").Value = 22270.65
$ws.Range("M31").Value = -21975.65
$ws.Range("H34").Value = 11931.978
$ws.Range("I34").Value = 22270.65
$ws.Range("K34").Value = 22270.65
$ws.Range("M34").Value = -22068.65
$ws.Range("H58").Value = 16577.594
$ws.Range("I58").Value = 974.75
$ws.Range("J58").Value = 63386.125
$ws.Range("K58").Value = 974.75
$ws.Range("L58").Value = 63386.125
$ws.Range("M58").Value = -771.75
$ws.Range("N58").Value = -63792.125
$ws.Range("H108").Value = 39669.75
$ws.Range("J108").Value = 39669.75
$ws.Range("L108").Value = 39669.75
$ws.Range("N108").Value = -47349.75
$ws.Range("H136").Value = 16577.594
$ws.Range("I136").Value = 974.75
$ws.Range("J136").Value = 63386.125
$ws.Range("K136").Value = 2924.25
$ws.Range("L136").Value = 190158.375
$ws.Range("M136").Value = -374.25
$ws.Range("N136").Value = -195258.375

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1834.6666
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 2252
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 6756
$ws.Range("M46").Value = -2909
$ws.Range("N46").Value = -6938
$ws.Range("H58").Value = 3567.5625
$ws.Range("I58").Value = 1499.5
$ws.Range("J58").Value = 3863
$ws.Range("K58").Value = 4498.5
$ws.Range("L58").Value = 11589
$ws.Range("M58").Value = -4370.5
$ws.Range("N58").Value = -11845
$ws.Range("H68").Value = 3695.3076
$ws.Range("J68").Value = 4023.3713
$ws.Range("L68").Value = 12070.1139
$ws.Range("N68").Value = -13692.1139
$ws.Range("H69").Value = 2500
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("H71").Value = 3695.3076
$ws.Range("J71").Value = 4023.3713
$ws.Range("L71").Value = 36210.3417
$ws.Range("N71").Value = -44322.3417
$ws.Range("H72").Value = 2500
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("H131").Value = 114467.01
$ws.Range("J131").Value = 125838.836
$ws.Range("L131").Value = 377516.508
$ws.Range("N131").Value = -387596.508

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 19400
$ws.Range("J15").Value = 19400
$ws.Range("L15").Value = 19400
$ws.Range("N15").Value = -19976
$ws.Range("H81").Value = 19400
$ws.Range("J81").Value = 19400
$ws.Range("L81").Value = 19400
$ws.Range("N81").Value = -21396
$ws.Range("H84").Value = 19400
$ws.Range("J84").Value = 19400
$ws.Range("L84").Value = 58200
$ws.Range("N84").Value = -68184
$ws.Range("H107").Value = 888.82355
$ws.Range("I107").Value = 378.7143
$ws.Range("J107").Value = 1245.9
$ws.Range("K107").Value = 378.7143
$ws.Range("L107").Value = 1245.9
$ws.Range("M107").Value = 1541.2857
$ws.Range("N107").Value = -5085.9

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1825.2858
$ws.Range("I22").Value = 1982.5
$ws.Range("K22").Value = 1982.5
$ws.Range("M22").Value = -1687.5
$ws.Range("H27").Value = 1825.2858
$ws.Range("I27").Value = 1982.5
$ws.Range("K27").Value = 1982.5
$ws.Range("M27").Value = -1875.5
$ws.Range("H132").Value = 1561.0889
$ws.Range("I132").Value = 1165.3636
$ws.Range("K132").Value = 3496.0908
$ws.Range("M132").Value = -966.0907999999999

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4455
$ws.Range("I62").Value = 2745
$ws.Range("J62").Value = 5432.143
$ws.Range("K62").Value = 2745
$ws.Range("L62").Value = 5432.143
$ws.Range("M62").Value = -2121
$ws.Range("N62").Value = -6680.143
$ws.Range("H65").Value = 4455
$ws.Range("I65").Value = 2745
$ws.Range("J65").Value = 5432.143
$ws.Range("K65").Value = 13725
$ws.Range("L65").Value = 27160.715
$ws.Range("M65").Value = -10605
$ws.Range("N65").Value = -33400.715
$ws.Range("H81").Value = 2366.8333
$ws.Range("I81").Value = 2440.2
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 4880.4
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -3819.4
$ws.Range("N81").Value = -6122
$ws.Range("H84").Value = 2366.8333
$ws.Range("I84").Value = 2440.2
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 24402
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -19098
$ws.Range("N84").Value = -30608
$ws.Range("H86").Value = 10000
$ws.Range("J86").Value = 10000
$ws.Range("L86").Value = 10000
$ws.Range("N86").Value = -12246
$ws.Range("H89").Value = 10000
$ws.Range("J89").Value = 10000
$ws.Range("L89").Value = 50000
$ws.Range("N89").Value = -61232
$ws.Range("H132").Value = 1720
$ws.Range("I132").Value = 1537.1111
$ws.Range("J132").Value = 2337.25
$ws.Range("K132").Value = 4611.3333
$ws.Range("L132").Value = 7011.75
$ws.Range("M132").Value = -2081.3333
$ws.Range("N132").Value = -12071.75

